$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 265
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
# Row 17
$ws.Range("H17").Value = 194584.48
$ws.Range("J17").Value = 194584.48
$ws.Range("L17").Value = 583753.4400000001
$ws.Range("N17").Value = -584089.4400000001
# Row 28
$ws.Range("H28").Value = 1778.2667
$ws.Range("I28").Value = 1063.1818
$ws.Range("J28").Value = 3744.75
$ws.Range("K28").Value = 1063.1818
$ws.Range("L28").Value = 3744.75
$ws.Range("M28").Value = -578.1818000000001
$ws.Range("N28").Value = -4714.75
# Row 32
$ws.Range("H32").Value = 862
$ws.Range("J32").Value = 843
$ws.Range("L32").Value = 843
$ws.Range("N32").Value = -1495
# Row 51
$ws.Range("H51").Value = 4082.5217
$ws.Range("I51").Value = 1471.1428
$ws.Range("J51").Value = 5225
$ws.Range("K51").Value = 1471.1428
$ws.Range("L51").Value = 5225
$ws.Range("M51").Value = -987.1428000000001
$ws.Range("N51").Value = -6193
# Row 62
$ws.Range("H62").Value = 4886.0625
$ws.Range("I62").Value = 3523.125
$ws.Range("K62").Value = 3523.125
$ws.Range("M62").Value = -2899.125
# Row 64
$ws.Range("H64").Value = 3875.5557
$ws.Range("J64").Value = 3910
$ws.Range("L64").Value = 3910
$ws.Range("N64").Value = -4406
# Row 65
$ws.Range("H65").Value = 4886.0625
$ws.Range("I65").Value = 3523.125
$ws.Range("K65").Value = 17615.625
$ws.Range("M65").Value = -14495.625
# Row 67
$ws.Range("H67").Value = 3875.5557
$ws.Range("J67").Value = 3910
$ws.Range("L67").Value = 3910
$ws.Range("N67").Value = -5626
# Row 86
$ws.Range("H86").Value = 5272.8667
$ws.Range("I86").Value = 5044.143
$ws.Range("J86").Value = 5473
$ws.Range("K86").Value = 5044.143
$ws.Range("L86").Value = 5473
$ws.Range("M86").Value = -3921.143
$ws.Range("N86").Value = -7719
# Row 89
$ws.Range("H89").Value = 5272.8667
$ws.Range("I89").Value = 5044.143
$ws.Range("J89").Value = 5473
$ws.Range("K89").Value = 25220.715
$ws.Range("L89").Value = 27365
$ws.Range("M89").Value = -19604.715
$ws.Range("N89").Value = -38597
# Row 98
$ws.Range("H98").Value = 1491.9286
$ws.Range("I98").Value = 1104.9231
$ws.Range("J98").Value = 1827.3334
$ws.Range("K98").Value = 1104.9231
$ws.Range("L98").Value = 1827.3334
$ws.Range("M98").Value = 393.0769
$ws.Range("N98").Value = -4823.3334
# Row 103
$ws.Range("H103").Value = 6665.5557
$ws.Range("I103").Value = 7387.5
$ws.Range("J103").Value = 890
$ws.Range("K103").Value = 22162.5
$ws.Range("L103").Value = 2670
$ws.Range("M103").Value = -21576.5
$ws.Range("N103").Value = -3842
# Row 113
$ws.Range("H113").Value = 4658.75
$ws.Range("I113").Value = 3599.8
$ws.Range("K113").Value = 3599.8
$ws.Range("M113").Value = -345.8000000000002
# Row 122
$ws.Range("H122").Value = 1491.9286
$ws.Range("I122").Value = 1104.9231
$ws.Range("J122").Value = 1827.3334
$ws.Range("K122").Value = 3314.7693
$ws.Range("L122").Value = 5482.0002
$ws.Range("M122").Value = -864.7692999999999
$ws.Range("N122").Value = -10382.0002
# Row 132
$ws.Range("H132").Value = 17095820
$ws.Range("I132").Value = 2137624.5
$ws.Range("J132").Value = 76928600
$ws.Range("K132").Value = 6412873.5
$ws.Range("L132").Value = 230785800
$ws.Range("M132").Value = -6410343.5
$ws.Range("N132").Value = -230790860

$ws = $wb.Worksheets.Item("ARM")
# Row 38
$ws.Range("H38").Value = 3393.4
$ws.Range("I38").Value = 3393.4
$ws.Range("K38").Value = 3393.4
$ws.Range("M38").Value = -2926.4

$ws = $wb.Worksheets.Item("BSM")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20001532
$ws.Range("I31").Value = 27778700
$ws.Range("K31").Value = 27778700
$ws.Range("M31").Value = -27778405
# Row 34
$ws.Range("H34").Value = 20001532
$ws.Range("I34").Value = 27778700
$ws.Range("K34").Value = 27778700
$ws.Range("M34").Value = -27778498
# Row 35
$ws.Range("H35").Value = 3858.3333
$ws.Range("I35").Value = 830
$ws.Range("J35").Value = 19000
$ws.Range("K35").Value = 830
$ws.Range("L35").Value = 19000
$ws.Range("M35").Value = -536
$ws.Range("N35").Value = -19588
# Row 62
$ws.Range("H62").Value = 2839.2856
$ws.Range("I62").Value = 2350
$ws.Range("J62").Value = 3328.5715
$ws.Range("K62").Value = 2350
$ws.Range("L62").Value = 3328.5715
$ws.Range("M62").Value = -1726
$ws.Range("N62").Value = -4576.5715
# Row 65
$ws.Range("H65").Value = 2839.2856
$ws.Range("I65").Value = 2350
$ws.Range("J65").Value = 3328.5715
$ws.Range("K65").Value = 11750
$ws.Range("L65").Value = 16642.8575
$ws.Range("M65").Value = -8630
$ws.Range("N65").Value = -22882.8575
# Row 94
$ws.Range("H94").Value = 1053.8889
$ws.Range("J94").Value = 855.1429000000001
$ws.Range("L94").Value = 855.1429000000001
$ws.Range("N94").Value = -1757.1429
# Row 99
$ws.Range("H99").Value = 2664.8604
$ws.Range("I99").Value = 2455.7856
$ws.Range("J99").Value = 3055.1333
$ws.Range("K99").Value = 2455.7856
$ws.Range("L99").Value = 3055.1333
$ws.Range("M99").Value = -957.7856000000002
$ws.Range("N99").Value = -6051.1333
# Row 119
$ws.Range("H119").Value = 31142.857
$ws.Range("J119").Value = 31142.857
$ws.Range("L119").Value = 31142.857
$ws.Range("N119").Value = -40818.857
# Row 122
$ws.Range("H122").Value = 1514.4736
$ws.Range("J122").Value = 1646.7142
$ws.Range("L122").Value = 4940.142599999999
$ws.Range("N122").Value = -9840.142599999999
# Row 126
$ws.Range("H126").Value = 2664.8604
$ws.Range("I126").Value = 2455.7856
$ws.Range("J126").Value = 3055.1333
$ws.Range("K126").Value = 7367.3568
$ws.Range("L126").Value = 9165.3999
$ws.Range("M126").Value = -4897.3568
$ws.Range("N126").Value = -14105.3999

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 1215.375
$ws.Range("I121").Value = 194
$ws.Range("J121").Value = 1361.2858
$ws.Range("K121").Value = 582
$ws.Range("L121").Value = 4083.8574
$ws.Range("M121").Value = 728
$ws.Range("N121").Value = -6703.857400000001
# Row 131
$ws.Range("H131").Value = 870.62
$ws.Range("J131").Value = 909.80896
$ws.Range("L131").Value = 2729.42688
$ws.Range("N131").Value = -12809.42688
# Row 137
$ws.Range("H137").Value = 23338164
$ws.Range("I137").Value = 2299.0908
$ws.Range("J137").Value = 37598972
$ws.Range("K137").Value = 6897.2724
$ws.Range("L137").Value = 112796916
$ws.Range("M137").Value = -1797.2724
$ws.Range("N137").Value = -112807116

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 4875
$ws.Range("I24").Value = 2000
$ws.Range("J24").Value = 8900
$ws.Range("K24").Value = 2000
$ws.Range("L24").Value = 8900
$ws.Range("M24").Value = -1827
$ws.Range("N24").Value = -9246
# Row 102
$ws.Range("H102").Value = 1225.3871
$ws.Range("I102").Value = 974.125
$ws.Range("J102").Value = 2086.8572
$ws.Range("K102").Value = 974.125
$ws.Range("L102").Value = 2086.8572
$ws.Range("M102").Value = 647.875
$ws.Range("N102").Value = -5330.8572
# Row 120
$ws.Range("H120").Value = 20000
$ws.Range("J120").Value = 20000
$ws.Range("L120").Value = 20000
$ws.Range("N120").Value = -29676

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1988.8889
$ws.Range("I100").Value = 2050
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 2050
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -1509
$ws.Range("N100").Value = -2582
# Row 121
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 330.4375
$ws.Range("I113").Value = 235.39131
$ws.Range("J113").Value = 573.3333
$ws.Range("K113").Value = 706.17393
$ws.Range("L113").Value = 1719.9999
$ws.Range("M113").Value = 1463.82607
$ws.Range("N113").Value = -6059.9999
# Row 121
$ws.Range("H121").Value = 30013.334
$ws.Range("J121").Value = 30013.334
$ws.Range("L121").Value = 30013.334
$ws.Range("N121").Value = -33507.334
# Row 136
$ws.Range("H136").Value = 18183656
$ws.Range("I136").Value = 23257622
$ws.Range("K136").Value = 69772866
$ws.Range("M136").Value = -69770316
